# Fixed columns length issue
#
# The source workbook carried a stray, empty drawing part
# (xl/drawings/drawing1.xml, an <xdr:wsDr/> with no shapes) left over from
# whatever tool produced the file. Real Excel does not expose a way to drop
# an orphaned drawing relationship through the object model when there are
# no Shapes to delete, so we rebuild the sheet: add a fresh sheet, copy the
# real data + header formatting across, drop the old sheet (and its
# drawing relationship with it), then rename the new sheet back to Sheet1.
$wb = $excel.ActiveWorkbook
$oldName = $wb.ActiveSheet.Name()

$newSheet = $wb.Worksheets.Add()
$oldSheet = $wb.Sheets.Item($oldName)

# Bring over the real table (A1:C4) and the styled-but-empty header cells
# (D1:Z1) so the rebuilt sheet matches the original layout exactly.
$oldSheet.Range("A1:C4").Copy($newSheet.Range("A1"))
$oldSheet.Range("D1:Z1").Copy($newSheet.Range("D1"))

[void]$wb.Sheets.Item($oldName).Delete()

$ws = $wb.Sheets.Item(1)
$ws.Name = $oldName
$ws.Activate()

# Add the new "Date" column.
$ws.Range("D1").Value = "Date"

$ws.Range("D2").Value = Get-Date -Year 1999 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Range("D3").Value = Get-Date -Year 2020 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("D4").Value = Get-Date -Year 2022 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0

[void]($ws.Range("D2:D4").NumberFormat = "m/d/yyyy")

[void]$ws.Range("D4").Select()
